# Rename unclear "dependents" month-style column headers (JAN..DEC)
# on the "Pool" worksheet to dependents_1..dependents_12.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

$names = @(
    "dependents_1",
    "dependents_2",
    "dependents_3",
    "dependents_4",
    "dependents_5",
    "dependents_6",
    "dependents_7",
    "dependents_8",
    "dependents_9",
    "dependents_10",
    "dependents_11",
    "dependents_12"
)

# Columns N (14) through Y (25) on row 1 hold the old JAN..DEC headers.
$startCol = 14
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $names[$i]
}
